$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "87.133.13"
$ws.Range("E2").Value = "  +10.03%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.356.38"
$ws.Range("E3").Value = "  +6.99%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - Solana
Set-TextValue "D5" "220.33"
$ws.Range("E5").Value = "  +8.09%  "

# Row 6 - BNB
Set-TextValue "D6" "638.69"
$ws.Range("E6").Value = "  +2.74%  "

# Row 7 - Dogecoin
Set-TextValue "D7" "0.328"
$ws.Range("E7").Value = "  +25.04%  "

# Row 8 - USDC
Set-TextValue "D8" "0.999"
$ws.Range("E8").Value = "  -0.12%  "

# Row 9 - XRP
Set-TextValue "D9" "0.625"
$ws.Range("E9").Value = "  +6.96%  "

# Row 10 - LidoStakedEther
Set-TextValue "D10" "3.391.39"
$ws.Range("E10").Value = "  +8.01%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.612"
$ws.Range("E11").Value = "  +6.04%  "

# Row 12 - ShibaInu
Set-TextValue "D12" "0.0000277"
$ws.Range("E12").Value = "  +11.64%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +2.14%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "3.973.60"
$ws.Range("E14").Value = "  +6.76%  "

# Row 15 - Avalanche
Set-TextValue "D15" "34.40"
$ws.Range("E15").Value = "  +11.17%  "

# Row 16 - Toncoin
$ws.Range("E16").Value = "  +4.03%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "87.107.30"
$ws.Range("E17").Value = "  +9.82%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "3.358.25"
$ws.Range("E18").Value = "  +6.84%  "

# Row 19 - Chainlink
Set-TextValue "D19" "14.77"

# Row 20 - SuiNetwork
Set-TextValue "D20" "3.23"
$ws.Range("E20").Value = "  +11.81%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "449.85"
$ws.Range("E21").Value = "  +4.01%  "

# Row 22 - Uniswap
Set-TextValue "D22" "9.21"
$ws.Range("E22").Value = "  +2.45%  "

# Row 23 - Polkadot
Set-TextValue "D23" "5.34"
$ws.Range("E23").Value = "  +3.29%  "

# Row 24 - LEO
Set-TextValue "D24" "7.42"
$ws.Range("E24").Value = "  +9.17%  "

# Row 25 - NEARProtocol
Set-TextValue "D25" "5.42"
$ws.Range("E25").Value = "  +17.46%  "

# Row 26 - Aptos
Set-TextValue "D26" "12.42"
$ws.Range("E26").Value = "  +16.15%  "

# Row 27 - WrappedeETH
Set-TextValue "D27" "3.457.59"
$ws.Range("E27").Value = "  +4.32%  "

# Row 28 - Litecoin
Set-TextValue "D28" "78.88"
$ws.Range("E28").Value = "  +4.58%  "

# Row 29 - PEPE
Set-TextValue "D29" "0.0000134"
$ws.Range("E29").Value = "  +12.16%  "

# Row 30 - Dai
$ws.Range("E30").Value = "  +0.03%  "

# Row 31 - Cronos
Set-TextValue "D31" "0.188"
$ws.Range("E31").Value = "  +54.97%  "

# Row 32 - Bittensor
Set-TextValue "D32" "607.33"
$ws.Range("E32").Value = "  +12.11%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +5.39%  "

# Rows 34/35 swap: Fetch.AI <-> Binance-PegBSC-USD
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D34" "1.00"
$ws.Range("E34").Value = "  +0.12%  "

$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D35" "1.57"
$ws.Range("E35").Value = "  +7.71%  "

# Row 36 - PancakeSwap
$ws.Range("E36").Value = "  +4.48%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  +3.15%  "

# Row 38 - EthereumClassic
Set-TextValue "D38" "23.62"

# Row 39 - RenderToken
Set-TextValue "D39" "6.58"
$ws.Range("E39").Value = "  +18.73%  "

# Row 40 - PolygonEcosystemToken
$ws.Range("E40").Value = "  +6.42%  "

# Row 41 - FirstDigitalUSD
Set-TextValue "D41" "1.00"
$ws.Range("E41").Value = "  +0.17%  "

# Rows 42/43 swap: dogwifhat <-> Stacks
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D42" "2.09"
$ws.Range("E42").Value = "  +17.61%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D43" "3.14"
$ws.Range("E43").Value = "  +18.73%  "

# Row 44 - WhiteBITCoin
Set-TextValue "D44" "21.35"
$ws.Range("E44").Value = "  +3.22%  "

# Row 45 - USDe
$ws.Range("E45").Value = "  +0.07%  "

# Rows 46/47 swap: Aave <-> Monero
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D46" "157.66"
$ws.Range("E46").Value = "  -3.10%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D47" "191.25"
$ws.Range("E47").Value = "  +2.83%  "

# Row 48 - ImmutableX
Set-TextValue "D48" "1.38"
$ws.Range("E48").Value = "  +7.62%  "

# Row 49 - OKB
Set-TextValue "D49" "45.83"
$ws.Range("E49").Value = "  +6.77%  "

# Row 50 - Mantle
Set-TextValue "D50" "0.794"
$ws.Range("E50").Value = "  +2.63%  "

# Row 51 - ARBITRUM
$ws.Range("E51").Value = "  +7.03%  "
